$wb = $excel.ActiveWorkbook

# Updated price/profit figures refreshed by the scheduled pricing run.
# Each block sets the changed cells for a given sheet/row based on the
# latest currentAveragePrice* / LevePrice* / LeveProfit* figures.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 7992
$ws.Range("I69").Value = 7992
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 23976
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -23102
$ws.Range("N69").ClearContents()  # was -211793

$ws.Range("H72").Value = 7992
$ws.Range("I72").Value = 7992
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 71928
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -67560
$ws.Range("N72").ClearContents()  # was -638871

$ws.Range("H112").Value = 1399.2
$ws.Range("I112").Value = 1000
$ws.Range("J112").Value = 1499
$ws.Range("K112").Value = 3000
$ws.Range("L112").Value = 4497
$ws.Range("M112").Value = -1892
$ws.Range("N112").Value = -6713

$ws.Range("H135").Value = 6436.263
$ws.Range("I135").Value = 1238.2778
$ws.Range("J135").Value = 100000
$ws.Range("K135").Value = 11144.5002
$ws.Range("L135").Value = 900000
$ws.Range("M135").Value = -8609.5002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1259
$ws.Range("I2").Value = 1259
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1259
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1146

$ws.Range("H32").Value = 22733624
$ws.Range("I32").Value = 23814272
$ws.Range("J32").Value = 39998
$ws.Range("K32").Value = 23814272
$ws.Range("L32").Value = 39998
$ws.Range("M32").Value = -23813985
$ws.Range("N32").Value = -40572

$ws.Range("H74").Value = 22744408
$ws.Range("I74").Value = 50000804
$ws.Range("J74").Value = 30745.334
$ws.Range("K74").Value = 50000804
$ws.Range("L74").Value = 30745.334
$ws.Range("M74").Value = -49999930

$ws.Range("H77").Value = 22744408
$ws.Range("I77").Value = 50000804
$ws.Range("J77").Value = 30745.334
$ws.Range("K77").Value = 250004020
$ws.Range("L77").Value = 153726.67
$ws.Range("M77").Value = -249999652

$ws.Range("H116").Value = 1259
$ws.Range("I116").Value = 1259
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1259
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1035

$ws.Range("H132").Value = 5302.9536
$ws.Range("I132").Value = 2389.4546
$ws.Range("J132").Value = 14917.5
$ws.Range("K132").Value = 7168.3638
$ws.Range("L132").Value = 44752.5
$ws.Range("M132").Value = -4638.3638
$ws.Range("N132").Value = -49812.5

$ws.Range("H133").Value = 70000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 70000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 70000
$ws.Range("N133").Value = -75060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1259
$ws.Range("I3").Value = 1259
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1259
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1145

$ws.Range("H20").Value = 2916.9546
$ws.Range("I20").Value = 2897.875
$ws.Range("J20").Value = 2967.8333
$ws.Range("K20").Value = 2897.875
$ws.Range("L20").Value = 2967.8333
$ws.Range("M20").Value = -2650.875
$ws.Range("N20").Value = -3461.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1002
$ws.Range("I16").Value = 1002
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1002
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -715

$ws.Range("H31").Value = 536921.25
$ws.Range("I31").Value = 9270.704
$ws.Range("J31").Value = 1374954.5
$ws.Range("K31").Value = 9270.704
$ws.Range("L31").Value = 1374954.5
$ws.Range("M31").Value = -8975.704

$ws.Range("H34").Value = 536921.25
$ws.Range("I34").Value = 9270.704
$ws.Range("J34").Value = 1374954.5
$ws.Range("K34").Value = 9270.704
$ws.Range("L34").Value = 1374954.5
$ws.Range("M34").Value = -9068.704

$ws.Range("H58").Value = 1332.2
$ws.Range("I58").Value = 1262.4
$ws.Range("J58").Value = 1402
$ws.Range("K58").Value = 1262.4
$ws.Range("L58").Value = 1402
$ws.Range("M58").Value = -1059.4
$ws.Range("N58").Value = -1808

$ws.Range("H107").Value = 464.45456
$ws.Range("I107").Value = 331.57144
$ws.Range("J107").Value = 697
$ws.Range("K107").Value = 331.57144
$ws.Range("L107").Value = 697
$ws.Range("M107").Value = 1588.42856

$ws.Range("H113").Value = 1002
$ws.Range("I113").Value = 1002
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1002
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1168

$ws.Range("H122").Value = 1578.8
$ws.Range("I122").Value = 1564.6666
$ws.Range("J122").Value = 1600
$ws.Range("K122").Value = 4693.9998
$ws.Range("L122").Value = 4800
$ws.Range("M122").Value = -2243.9998

$ws.Range("H135").Value = 90866.664
$ws.Range("I135").Value = 50000
$ws.Range("J135").Value = 99040
$ws.Range("K135").Value = 50000
$ws.Range("L135").Value = 99040
$ws.Range("M135").Value = -44930
$ws.Range("N135").Value = -109180

$ws.Range("H136").Value = 1332.2
$ws.Range("I136").Value = 1262.4
$ws.Range("J136").Value = 1402
$ws.Range("K136").Value = 3787.2
$ws.Range("L136").Value = 4206
$ws.Range("M136").Value = -1237.2
$ws.Range("N136").Value = -9306

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 522
$ws.Range("I107").Value = 403.4
$ws.Range("J107").Value = 670.25
$ws.Range("K107").Value = 1210.2
$ws.Range("L107").Value = 2010.75
$ws.Range("M107").Value = 709.8000000000002

$ws.Range("H120").Value = 22851.834
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 22851.834
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 68555.50199999999
$ws.Range("M120").ClearContents()  # was -7162
$ws.Range("N120").Value = -78231.50199999999

$ws.Range("H134").Value = 3557.95
$ws.Range("I134").Value = 3557.95
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 10673.85
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -5603.849999999999
$ws.Range("N134").ClearContents()  # was -31015.8

$ws.Range("H136").Value = 6951
$ws.Range("I136").Value = 6951
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 20853
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -15753

$ws.Range("H137").Value = 5925
$ws.Range("I137").Value = 7266.6665
$ws.Range("J137").Value = 1900
$ws.Range("K137").Value = 21799.9995
$ws.Range("L137").Value = 5700
$ws.Range("M137").Value = -16699.9995
$ws.Range("N137").Value = -15900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1618.7858
$ws.Range("I122").Value = 1561.7273
$ws.Range("J122").Value = 1828
$ws.Range("K122").Value = 4685.1819
$ws.Range("L122").Value = 5484
$ws.Range("M122").Value = -2235.1819
$ws.Range("N122").Value = -10384

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2172
$ws.Range("I22").Value = 2863.5
$ws.Range("J22").Value = 1250
$ws.Range("K22").Value = 2863.5
$ws.Range("L22").Value = 1250
$ws.Range("M22").Value = -2568.5
$ws.Range("N22").Value = -1840

$ws.Range("H27").Value = 2172
$ws.Range("I27").Value = 2863.5
$ws.Range("J27").Value = 1250
$ws.Range("K27").Value = 2863.5
$ws.Range("L27").Value = 1250
$ws.Range("M27").Value = -2756.5
$ws.Range("N27").Value = -1464

$ws.Range("H46").Value = 6060.409
$ws.Range("I46").Value = 3705.7144
$ws.Range("J46").Value = 10181.125
$ws.Range("K46").Value = 3705.7144
$ws.Range("L46").Value = 10181.125
$ws.Range("M46").Value = -3517.7144
$ws.Range("N46").Value = -10557.125

$ws.Range("H48").Value = 18498.334
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 18498.334
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 18498.334
$ws.Range("N48").Value = -19820.334

$ws.Range("H132").Value = 75421.12
$ws.Range("I132").Value = 9082.429
$ws.Range("J132").Value = 385001.66
$ws.Range("K132").Value = 27247.287
$ws.Range("L132").Value = 1155004.98
$ws.Range("M132").Value = -24717.287
$ws.Range("N132").Value = -1160064.98

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 745
$ws.Range("I81").Value = 745
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1490
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -429

$ws.Range("H84").Value = 745
$ws.Range("I84").Value = 745
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 7450
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -2146

$ws.Range("H126").Value = 5918.4
$ws.Range("I126").Value = 4614
$ws.Range("J126").Value = 7875
$ws.Range("K126").Value = 13842
$ws.Range("L126").Value = 23625
$ws.Range("M126").Value = -11372
